$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 15 header: add center/middle alignment to the already-bordered header
# cells (font/border stay as-is, only alignment is added).
# ---------------------------------------------------------------------------
$ws.Range("F15").HorizontalAlignment = -4108
$ws.Range("F15").VerticalAlignment = -4108
$ws.Range("G15:L15").HorizontalAlignment = -4108
$ws.Range("G15:L15").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 22 area: rework the "x rata-rata" / "y rata-rata" scratch-work block.
# ---------------------------------------------------------------------------

# Row 23: drop E23, F23 becomes x-bar label, H23/I23 get a distinct (left/top)
# style and new formulas.
$ws.Range("E23").ClearContents()
$ws.Range("F23").Value = "x" + [char]0x0304
$ws.Range("G23").Value = "Σx"
$ws.Range("H23").Formula = "=G16 /""12"""
$ws.Range("I23").Formula = "=H23/H24"

$ws.Range("H23:I23").Font.Name = "Times New Roman"
$ws.Range("H23:I23").Font.Size = 12
$ws.Range("H23:I23").Font.Bold = $false
$ws.Range("H23:I23").HorizontalAlignment = -4131
$ws.Range("H23:I23").VerticalAlignment = -4160

# Row 24: n label stays, value stays, add I24 (blank, styled) and re-style H24.
$ws.Range("G24").Value = "n"
$ws.Range("H24").Value = 12
$ws.Range("H24:I24").Font.Name = "Times New Roman"
$ws.Range("H24:I24").Font.Size = 12
$ws.Range("H24:I24").Font.Bold = $false
$ws.Range("H24:I24").HorizontalAlignment = -4131
$ws.Range("H24:I24").VerticalAlignment = -4160

# Rows 25-26: new blank styled filler rows.
$ws.Range("H25:I26").Font.Name = "Times New Roman"
$ws.Range("H25:I26").Font.Size = 12
$ws.Range("H25:I26").Font.Bold = $false
$ws.Range("H25:I26").HorizontalAlignment = -4131
$ws.Range("H25:I26").VerticalAlignment = -4160

# Row 27: y-bar label, Σy label, formulas re-styled to match rows 23-26.
$ws.Range("F27").Value = [char]0x0233
$ws.Range("G27").Value = "Σy"
$ws.Range("H27").Formula = "=H16"
$ws.Range("I27").Formula = "=H27/H28"
$ws.Range("H27:I27").Font.Name = "Times New Roman"
$ws.Range("H27:I27").Font.Size = 12
$ws.Range("H27:I27").Font.Bold = $false
$ws.Range("H27:I27").HorizontalAlignment = -4131
$ws.Range("H27:I27").VerticalAlignment = -4160

# Row 28: n label, value, and new I28 blank styled cell.
$ws.Range("G28").Value = "n"
$ws.Range("H28").Value = 12
$ws.Range("H28:I28").Font.Name = "Times New Roman"
$ws.Range("H28:I28").Font.Size = 12
$ws.Range("H28:I28").Font.Bold = $false
$ws.Range("H28:I28").HorizontalAlignment = -4131
$ws.Range("H28:I28").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Rows 30-32: shifted down from the previous 30-32 text (same text, same
# relative order -- "Persamaan garis yang memenuhi" / "Y = a + bx" / "Dengan :").
# ---------------------------------------------------------------------------
$ws.Range("F30").Value = "Persamaan garis yang memenuhi"
$ws.Range("F31").Value = "Y = a + bx"
$ws.Range("F32").Value = "Dengan :"

# ---------------------------------------------------------------------------
# Row 33-34: b = ... block; K/L columns get the new (left-aligned, no border)
# style.
# ---------------------------------------------------------------------------
$ws.Range("F33").Value = "b ="
$ws.Range("K33").Formula = "=12*I16 - (G16*H16)"
$ws.Range("L33").Formula = "=K33/K34"
$ws.Range("K34").Formula = "=12*J16 - (G16)^2"

$ws.Range("K33:L34").Font.Name = "Times New Roman"
$ws.Range("K33:L34").Font.Size = 12
$ws.Range("K33:L34").Font.Bold = $false
$ws.Range("K33:L34").HorizontalAlignment = -4131

# Rows 35-37: new blank filler cells in K/L with the same style.
$ws.Range("K35:L37").Font.Name = "Times New Roman"
$ws.Range("K35:L37").Font.Size = 12
$ws.Range("K35:L37").Font.Bold = $false
$ws.Range("K35:L37").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# Row 36: a = ... ; recompute with the corrected formula result.
# ---------------------------------------------------------------------------
$ws.Range("G36").Value = [char]0x0233 + " -bx"
$ws.Range("H36").Formula = "=I27-L33*I23"

# ---------------------------------------------------------------------------
# Rows 38-42: conclusion text, re-sequenced with a new row 38 + a gap row 39.
# ---------------------------------------------------------------------------
$ws.Range("F38").Value = "Jadi, b = 0,389106675 dan a = 31,82148"
$ws.Range("F39").ClearContents()
$ws.Range("F40").Value = "Sehingga persamaan garisnya :"
$ws.Range("F41").Value = "y = a + bx"
$ws.Range("F42").Value = "y = 31,82148 + 0,389106675x"

# ---------------------------------------------------------------------------
# Sheet view tweaks.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 66
$ws.Range("E32").Select()
